$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column A (rows 2:5) as text so the date-looking string
# "07/07/2023" is stored literally instead of being parsed into a date serial.
$ws.Range("A2:A5").NumberFormat = "@"

# Row 2 - numeric values
$ws.Range("A2").Value = "07/07/2023"
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 1000
$ws.Range("D2").Value = 1000
$ws.Range("E2").Value = 1000
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 100

# Row 3 - numeric values
$ws.Range("A3").Value = "07/07/2023"
$ws.Range("B3").Value = 1000
$ws.Range("C3").Value = 2000
$ws.Range("D3").Value = 2000
$ws.Range("E3").Value = 3000
$ws.Range("F3").Value = 20
$ws.Range("G3").Value = 1000
$ws.Range("H3").Value = 150

# Row 4 - numeric values
$ws.Range("A4").Value = "07/07/2023"
$ws.Range("B4").Value = 1000
$ws.Range("C4").Value = 3000
$ws.Range("D4").Value = 3000
$ws.Range("E4").Value = 6000
$ws.Range("F4").Value = 40
$ws.Range("G4").Value = 3000
$ws.Range("H4").Value = 200

# Row 5 - text values (stored as strings that look like numbers)
$ws.Range("A5").Value = "07/07/2023"
$ws.Range("B5:H5").NumberFormat = "@"
$ws.Range("B5").Value = "5000.00"
$ws.Range("C5").Value = "8000.00"
$ws.Range("D5").Value = "1000.00"
$ws.Range("E5").Value = "7000.00"
$ws.Range("F5").Value = "50.0"
$ws.Range("G5").Value = "1000.00"
$ws.Range("H5").Value = "87.50"
